$d = $word.ActiveDocument

$d.Content.Find.Execute("19÷6=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "45÷6=7, 3", 2)
$d.Content.Find.Execute("88÷2=44, 0", $true, $false, $false, $false, $false, $true, 1, $false, "76÷5=15, 1", 2)
$d.Content.Find.Execute("39÷6=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "90÷9=10, 0", 2)
$d.Content.Find.Execute("61÷9=6, 7", $true, $false, $false, $false, $false, $true, 1, $false, "34÷7=4, 6", 2)
$d.Content.Find.Execute("68÷5=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "65÷6=10, 5", 2)
$d.Content.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "46÷9=5, 1", 2)
$d.Content.Find.Execute("33÷2=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "90÷2=45, 0", 2)
$d.Content.Find.Execute("18÷6=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=14, 2", 2)
$d.Content.Find.Execute("31÷2=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "67÷3=22, 1", 2)
$d.Content.Find.Execute("74÷9=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "47÷5=9, 2", 2)
$d.Content.Find.Execute("68÷8=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "28÷5=5, 3", 2)
$d.Content.Find.Execute("33÷3=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "95÷7=13, 4", 2)
$d.Content.Find.Execute("51÷4=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "28÷2=14, 0", 2)
$d.Content.Find.Execute("58÷5=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "85÷2=42, 1", 2)
$d.Content.Find.Execute("13÷2=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷5=8, 0", 2)
$d.Content.Find.Execute("87÷4=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 2)
$d.Content.Find.Execute("53÷9=5, 8", $true, $false, $false, $false, $false, $true, 1, $false, "19÷8=2, 3", 2)
$d.Content.Find.Execute("82÷2=41, 0", $true, $false, $false, $false, $false, $true, 1, $false, "37÷5=7, 2", 2)
$d.Content.Find.Execute("18÷9=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷7=12, 5", 2)
$d.Content.Find.Execute("45÷7=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "81÷8=10, 1", 2)
$d.Content.Find.Execute("96÷8=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=4, 4", 2)
$d.Content.Find.Execute("19÷7=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "95÷9=10, 5", 2)
$d.Content.Find.Execute("71÷3=23, 2", $true, $false, $false, $false, $false, $true, 1, $false, "85÷3=28, 1", 2)
$d.Content.Find.Execute("88÷5=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=7, 1", 2)
$d.Content.Find.Execute("41÷8=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=7, 0", 2)
